$d = $word.ActiveDocument

# Change 1: "a lower error than the remaining points." -> "lower errors than the remaining points."
[void]$d.Content.Find.Execute("a lower error than the remaining points.", $true, $false, $false, $false, $false,
                               $true, 1, $false, "lower errors than the remaining points.", 2)

# Change 2: " at the end of the test for data quality concerns, " -> " at the end of testing due to data quality concerns, "
[void]$d.Content.Find.Execute(" at the end of the test for data quality concerns, ", $true, $false, $false, $false, $false,
                               $true, 1, $false, " at the end of testing due to data quality concerns, ", 2)
